$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, entered in the same order the original author created the
# --- shared strings table (so sharedStrings.xml indices line up exactly) ---
$ws.Range("B3").Value = "Backlog:"
$ws.Range("C3").Value = "User Storys:"
$ws.Range("D3").Value = "To do:"
$ws.Range("E3").Value = "Doing:"
$ws.Range("F3").Value = "Review:"
$ws.Range("G3").Value = "Done:"

$ws.Range("E4").Value = "Dar ideias no servidor de discord e discuti-las"
$ws.Range("G7").Value = "Fazer a hierarquia de ficheiros e adiconá-los"
$ws.Range("G5").Value = "Fazer fork ao projeto e adicionar os membros ao mesmo"
$ws.Range("G6").Value = "Meeting semanal"
$ws.Range("G4").Value = "Fazer servidor de discord para o trabalho e organiza-lo"
$ws.Range("D5").Value = "Analisar o código dado"
$ws.Range("D4").Value = "Começar a fazer o use case diagram"
$ws.Range("E5").Value = "Jogar o jogo para conhecer melhor o projeto"
$ws.Range("F4").Value = "Alterar o readme file no repositório git"
$ws.Range("B4").Value = "Como um jogador novo deste tipo de jogos eu quero um tutorial para poder perceber como se começa a jogar"
$ws.Range("B5").Value = "Como um jogador já com alguma experiência e horas no jogo gostava que houvessem casamentos e estes formacem alianças para expandir a paz e comércio"
$ws.Range("B6").Value = "Como um jogador já com alguma experiência gostava de ter um sistema de preços conforme a demanda para que não possa abusar de loopholes"

# --- Header row styling: Excel's built-in "Note" cell style (yellow fill +
# --- thin gray border), with a bold dark-gray font override ---
$hdr = $ws.Range("B3:G3")
$hdr.Style = "Note"
$hdr.Font.Bold = $true
$hdr.Font.Size = 11
$hdr.Font.Color = 4144959

# --- B4 just carries a (no-op) alignment flag in the source file ---
$ws.Range("B4").HorizontalAlignment = 1

# --- Column widths ---
$ws.Columns("B").ColumnWidth = 139
$ws.Columns("C").ColumnWidth = 22.6666666
$ws.Columns("D").ColumnWidth = 40
$ws.Columns("E").ColumnWidth = 39.8333333
$ws.Columns("F").ColumnWidth = 36.5
$ws.Columns("G").ColumnWidth = 51

# --- Row 4 height ---
$ws.Rows(4).RowHeight = 15.75

# --- Selection state ---
$ws.Range("B4").Select()

Write-Output "done"
